$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A189").EntireRow.Insert()

$ws.Range("A189").Value = 3
$ws.Range("B189").Value = "Femacal de La Calera"
$ws.Range("C189").Value = "Coquimbo"
$ws.Range("D189").Value = 44504
$ws.Range("E189").Value = 5
$ws.Range("F189").Value = 100112009
$ws.Range("G189").Value = "Acelga"
$ws.Range("H189").Value = "Sin especificar"
$ws.Range("I189").Value = "Primera"
$ws.Range("J189").Value = 230
$ws.Range("K189").Value = 2000
$ws.Range("L189").Value = 2200
$ws.Range("M189").Value = 2096
$ws.Range("N189").Value = '$/docena de atados (6 kilos)'
$ws.Range("O189").Value = "Provincia de Quillota"
$ws.Range("P189").Value = 349
$ws.Range("Q189").Value = 6
$ws.Range("R189").Value = "Hortaliza"
